# Shrink the font size of the four "stage_*" label text boxes on every
# slide from 12pt to 11pt, and re-fit the (spAutoFit) text box height to
# match the smaller font, matching how PowerPoint lays out an
# auto-sized text box after an in-place font-size edit.

$p = $ppt.ActivePresentation

# 13.3289 pt == 169277 EMU once PowerPoint round-trips the Height
# property through its single-precision internal representation.
$newHeightPt = 13.3289

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)

        if ($shape.Name -like "TextBox*") {
            $textRange = $shape.TextFrame.TextRange
            $textRange.Font.Size = 11

            $shape.Height = $newHeightPt
        }
    }
}
